$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "MONTO MMOO" column header in G1 -------------------------
# Copy the header formatting from F1 (same font/alignment as the other
# header cells) onto G1, then overwrite the value so a fresh shared string
# "MONTO MMOO" is created.
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("G1").Value = "MONTO MMOO"

# --- Drop the stray border-only formatting that used to cover E2:G10 ------
# Those cells carried no data, only a blank "bordered" style; clearing them
# removes the cell records entirely (matching the now-unformatted blanks).
$ws.Range("E2:G10").Clear()

# --- Remove the trailing all-blank rows 11-13 ------------------------------
$ws.Rows("11:13").Delete()

# --- Mark the new MMOO entry field (G6) with an underlined style ----------
$ws.Range("G6").Font.Underline = 2

# --- Give the new column a sensible width ----------------------------------
$ws.Columns("G").ColumnWidth = 13.83

# --- Move the active selection to the new field ----------------------------
$ws.Range("G6").Select() | Out-Null
